$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 364, shifting existing rows 364-436 down to 365-437.
$ws.Rows(364).Insert()

# Populate the new row 364 with the new data record.
$ws.Range("A364").Value = 4
$ws.Range("B364").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C364").Value = "Los Lagos"
$ws.Range("D364").Value = 45173
$ws.Range("E364").Value = 10
$ws.Range("F364").Value = "Fruta"
$ws.Range("G364").Value = 100108
$ws.Range("H364").Value = "Tropicales y subtropicales"
$ws.Range("I364").Value = 100108002
$ws.Range("J364").Value = "Mango"
$ws.Range("K364").Value = "Sin especificar"
$ws.Range("L364").Value = "Primera"
$ws.Range("M364").Value = 40
$ws.Range("N364").Value = 12000
$ws.Range("O364").Value = 12000
$ws.Range("P364").Value = 12000
$ws.Range("Q364").Value = "$/bandeja 4 kilos"
$ws.Range("R364").Value = "Brasil"
$ws.Range("S364").Value = 3000
$ws.Range("T364").Value = 4
